$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22, pushing the "(-) Front-End" checklist block down
$ws.Rows("22:22").Insert()

# Resize the table (ListObject) to include the newly inserted row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:G25"))

# Row 20: mark status as "Ok"
$ws.Range("B20").Value = "Ok"

# New row 22: "Versão Produção (Arrumar Inconsistencias)"
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = "Versão Produção (Arrumar Inconsistencias)"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("G23").Select() | Out-Null
